$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in column D (on_depot)
$ws.Range("D2").Value = 8
$ws.Range("D3").Value = 4
$ws.Range("D7").Value = 10

# Update selection to D2
$ws.Range("D2").Select()
